$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -7
$ws.Range("F8").Value = -4
$ws.Range("F10").Value = -5
$ws.Range("F12").Value = -11
$ws.Range("F14").Value = -7
$ws.Range("F16").Value = -6
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = -8
$ws.Range("F20").Value = -2
$ws.Range("F22").Value = 5
$ws.Range("F23").Value = -2
$ws.Range("F24").Value = -1
$ws.Range("F29").Value = -3
$ws.Range("F32").Value = 3
$ws.Range("F33").Value = -8
$ws.Range("F37").Value = -6
$ws.Range("F38").Value = 7
$ws.Range("F40").Value = -2
$ws.Range("F49").Value = 7
$ws.Range("F51").Value = 10
$ws.Range("F52").Value = 0
$ws.Range("F54").Value = 3
$ws.Range("F55").Value = -7
$ws.Range("F56").Value = -3
$ws.Range("F57").Value = -5
